$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed date) column C for rows 2-10 from 45221 to 45224
$ws.Range("C2:C10").Value = 45224
